# Add a new day's worth of hospital-occupancy data (2020-05-29, serial 43980)
# by duplicating the previous day's 20-row block (rows 1093-1112) into the
# next 20 rows (1113-1132), then updating the date and the bed-count values
# that changed for the new day.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Duplicate the last day's block (rows 1093-1112) into rows 1113-1132.
# This carries over hospital name, municipio, provincia, codigo_ine,
# observaciones values and all cell styles/number formats intact.
$src = $ws.Range("A1093:H1112")
$dst = $ws.Range("A1113:H1132")
$src.Copy($dst)

# Update the date column for the new block to the new day (2020-05-29).
$ws.Range("A1113:A1132").Value2 = 43980

# Update camas_ocupadas_total (C) / camas_uci_ocupadas (D) figures that
# changed for this new reporting day.
$ws.Cells.Item(1117, 3).ClearContents()   # Hospital General de la Defensa: now blank
$ws.Cells.Item(1118, 3).Value2 = 4        # Hospital Obispo Polanco
$ws.Cells.Item(1119, 3).Value2 = 1        # Hospital de Alcañiz
$ws.Cells.Item(1120, 3).Value2 = 11       # Hospital de Barbastro
$ws.Cells.Item(1121, 3).Value2 = 6        # Hospital San Jorge
$ws.Cells.Item(1124, 3).Value2 = 10       # Hospital San José
$ws.Cells.Item(1126, 3).Value2 = 3        # MAZ
